$d = $word.ActiveDocument

$replacements = @(
    @('151÷5=30, 1', '592÷8=74, 0'),
    @('684÷6=114, 0', '843÷4=210, 3'),
    @('140÷4=35, 0', '256÷8=32, 0'),
    @('119÷6=19, 5', '803÷9=89, 2'),
    @('741÷6=123, 3', '489÷4=122, 1'),
    @('789÷5=157, 4', '974÷6=162, 2'),
    @('979÷7=139, 6', '184÷2=92, 0'),
    @('152÷5=30, 2', '567÷3=189, 0'),
    @('571÷4=142, 3', '853÷3=284, 1'),
    @('708÷4=177, 0', '971÷8=121, 3'),
    @('940÷4=235, 0', '593÷9=65, 8'),
    @('664÷5=132, 4', '576÷4=144, 0'),
    @('300÷9=33, 3', '441÷8=55, 1'),
    @('757÷8=94, 5', '994÷2=497, 0'),
    @('870÷5=174, 0', '566÷8=70, 6'),
    @('301÷3=100, 1', '233÷4=58, 1'),
    @('796÷5=159, 1', '743÷8=92, 7'),
    @('711÷3=237, 0', '266÷2=133, 0'),
    @('373÷4=93, 1', '226÷2=113, 0'),
    @('402÷4=100, 2', '305÷8=38, 1'),
    @('450÷5=90, 0', '944÷5=188, 4'),
    @('218÷5=43, 3', '296÷4=74, 0'),
    @('966÷3=322, 0', '951÷3=317, 0'),
    @('527÷5=105, 2', '361÷5=72, 1'),
    @('772÷5=154, 2', '347÷2=173, 1'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

